$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old extra columns (X:AG) that no longer exist in the target ---
$ws.Range("X1:AG19").Clear()

# --- Row 1: numeric column headers B1:W1 (0..21), unchanged values, just ensure range ---
$row1 = New-Object 'object[,]' 1,22
$row1[0,0] = 0
$row1[0,1] = 1
$row1[0,2] = 2
$row1[0,3] = 3
$row1[0,4] = 4
$row1[0,5] = 5
$row1[0,6] = 6
$row1[0,7] = 7
$row1[0,8] = 8
$row1[0,9] = 9
$row1[0,10] = 10
$row1[0,11] = 11
$row1[0,12] = 12
$row1[0,13] = 13
$row1[0,14] = 14
$row1[0,15] = 15
$row1[0,16] = 16
$row1[0,17] = 17
$row1[0,18] = 18
$row1[0,19] = 19
$row1[0,20] = 20
$row1[0,21] = 21
$ws.Range("B1:W1").Value = $row1

# --- Row 2: text headers B2 (HKL) and C2:W2 (reflection / pair labels) ---
$row2 = New-Object 'object[,]' 1,22
$row2[0,0] = "HKL"
$row2[0,1] = '[4, 2, 0]'
$row2[0,2] = '[4, 0, 0]'
$row2[0,3] = '[2, 0, 0]'
$row2[0,4] = '[2, 2, 0]'
$row2[0,5] = '[3, 3, 3]'
$row2[0,6] = '[1, 1, 1]'
$row2[0,7] = '[2, 2, 2]'
$row2[0,8] = '[3, 3, 1]'
$row2[0,9] = '[3, 1, 1]'
$row2[0,10] = '[4, 2, 2]'
$row2[0,11] = '[5, 1, 1]'
$row2[0,12] = '1Pair-A'
$row2[0,13] = '1Pair-B'
$row2[0,14] = '2Pairs-A'
$row2[0,15] = '2Pairs-B'
$row2[0,16] = '3Pairs-A'
$row2[0,17] = '3Pairs-B'
$row2[0,18] = '3Pairs-C'
$row2[0,19] = '4Pairs'
$row2[0,20] = '5A4F'
$row2[0,21] = 'MaxUnique'
$ws.Range("B2:W2").Value = $row2

# --- Column A: row index numbers, rows 2..23 (0..21) ---
$colA = New-Object 'object[,]' 22,1
$colA[0,0] = 0
$colA[1,0] = 1
$colA[2,0] = 2
$colA[3,0] = 3
$colA[4,0] = 4
$colA[5,0] = 5
$colA[6,0] = 6
$colA[7,0] = 7
$colA[8,0] = 8
$colA[9,0] = 9
$colA[10,0] = 10
$colA[11,0] = 11
$colA[12,0] = 12
$colA[13,0] = 13
$colA[14,0] = 14
$colA[15,0] = 15
$colA[16,0] = 16
$colA[17,0] = 17
$colA[18,0] = 18
$colA[19,0] = 19
$colA[20,0] = 20
$colA[21,0] = 21
$ws.Range("A2:A23").Value = $colA

# --- Column B: scheme-name labels, rows 2..23 ---
$colB = New-Object 'object[,]' 22,1
$colB[0,0] = 'HKL'
$colB[1,0] = 'BT8Hex_2.5'
$colB[2,0] = 'BT8Hex_5'
$colB[3,0] = 'BT8Hex_10'
$colB[4,0] = 'BT8Hex_15'
$colB[5,0] = 'Spiral2.5'
$colB[6,0] = 'Spiral5'
$colB[7,0] = 'Spiral7.5'
$colB[8,0] = 'Spiral10'
$colB[9,0] = 'Spiral15'
$colB[10,0] = 'OffsetF45'
$colB[11,0] = 'OffsetA45'
$colB[12,0] = 'OffsetFTD'
$colB[13,0] = 'OffsetATD'
$colB[14,0] = 'Holden2.5'
$colB[15,0] = 'Holden5'
$colB[16,0] = 'Holden10'
$colB[17,0] = 'Holden15'
$colB[18,0] = 'HexGrid-90degTilt2.5degRes'
$colB[19,0] = 'HexGrid-90degTilt5degRes'
$colB[20,0] = 'HexGrid-90degTilt10degRes'
$colB[21,0] = 'HexGrid-90degTilt15degRes'
$ws.Range("B2:B23").Value = $colB

# --- Main data block C3:W23 ---
$data = New-Object 'object[,]' 21,21
$data[0,0] = 0.9995919374725618
$data[0,1] = 0.998417457385777
$data[0,2] = 0.998417457385777
$data[0,3] = 1.000252584110811
$data[0,4] = 1.001309339523482
$data[0,5] = 1.001309339523482
$data[0,6] = 1.001309339523482
$data[0,7] = 1.000572410175627
$data[0,8] = 0.9996513516244195
$data[0,9] = 1.000475110494009
$data[0,10] = 0.9989462551925379
$data[0,11] = 1.001309339523482
$data[0,12] = 1.000252584110811
$data[0,13] = 0.9993350207482943
$data[0,14] = 0.9999519678676154
$data[0,15] = 0.9999931270066901
$data[0,16] = 0.9994404643736693
$data[0,17] = 0.9999931270066901
$data[0,18] = 0.9999076831611224
$data[0,19] = 1.000188014433594
$data[0,20] = 0.9999020557474031
$data[1,0] = 0.9992105715851728
$data[1,1] = 0.99694737855699
$data[1,2] = 0.99694737855699
$data[1,3] = 1.000483617637495
$data[1,4] = 1.002532032879154
$data[1,5] = 1.002532032879154
$data[1,6] = 1.002532032879154
$data[1,7] = 1.001103762898346
$data[1,8] = 0.9993272647152831
$data[1,9] = 1.000918454697044
$data[1,10] = 0.9979667700354508
$data[1,11] = 1.002532032879154
$data[1,12] = 1.000483617637495
$data[1,13] = 0.9987154980972426
$data[1,14] = 0.9999054411763891
$data[1,15] = 0.9999876763578798
$data[1,16] = 0.9989194203032561
$data[1,17] = 0.9999876763578798
$data[1,18] = 0.9998225734472306
$data[1,19] = 1.000364465333615
$data[1,20] = 0.9998112316256169
$data[2,0] = 0.9985036944963006
$data[2,1] = 0.9941257153295334
$data[2,2] = 0.9941257153295334
$data[2,3] = 1.000966309005693
$data[2,4] = 1.004809098569759
$data[2,5] = 1.004809098569759
$data[2,6] = 1.004809098569759
$data[2,7] = 1.002127757514392
$data[2,8] = 0.9987075729724719
$data[2,9] = 1.001747605040423
$data[2,10] = 0.9960935419524207
$data[2,11] = 1.004809098569759
$data[2,12] = 1.000966309005693
$data[2,13] = 0.997546012167613
$data[2,14] = 0.9998369409890823
$data[2,15] = 0.9999670409683284
$data[2,16] = 0.9979331991025658
$data[2,17] = 0.9999670409683284
$data[2,18] = 0.9996521739693642
$data[2,19] = 1.000683558889443
$data[2,20] = 0.9996351618601242
$data[3,0] = 0.9978133812074885
$data[3,1] = 0.9913725499258402
$data[3,2] = 0.9913725499258402
$data[3,3] = 1.001436347624764
$data[3,4] = 1.007032548113844
$data[3,5] = 1.007032548113844
$data[3,6] = 1.007032548113844
$data[3,7] = 1.003126794342685
$data[3,8] = 0.9981028668106035
$data[3,9] = 1.002557141953514
$data[3,10] = 0.9942656314795321
$data[3,11] = 1.007032548113844
$data[3,12] = 1.001436347624764
$data[3,13] = 0.9964044487753023
$data[3,14] = 0.9997696072176838
$data[3,15] = 0.9999471485548161
$data[3,16] = 0.9969705881204027
$data[3,17] = 0.9999471485548161
$data[3,18] = 0.999486078118763
$data[3,19] = 1.000995372117779
$data[3,20] = 0.9994634076822839
$data[4,0] = 0.9999955526119443
$data[4,1] = 0.999941758293099
$data[4,2] = 0.999941758293099
$data[4,3] = 1.000025814220723
$data[4,4] = 1.000018818678323
$data[4,5] = 1.000018818678323
$data[4,6] = 1.000018818678323
$data[4,7] = 1.000022800664458
$data[4,8] = 0.9999881614183147
$data[4,9] = 1.000008310364364
$data[4,10] = 0.999964081434491
$data[4,11] = 1.000018818678323
$data[4,12] = 1.000025814220723
$data[4,13] = 0.9999837862569109
$data[4,14] = 1.000006987819519
$data[4,15] = 0.9999954637307148
$data[4,16] = 0.9999852446440455
$data[4,17] = 0.9999954637307148
$data[4,18] = 0.9999936381526148
$data[4,19] = 0.9999986742577563
$data[4,20] = 0.9999956622107147
$data[5,0] = 1.000008613809069
$data[5,1] = 0.9999629627462964
$data[5,2] = 0.9999629627462964
$data[5,3] = 1.000034294730554
$data[5,4] = 0.999980181993475
$data[5,5] = 0.999980181993475
$data[5,6] = 0.999980181993475
$data[5,7] = 1.000016376142092
$data[5,8] = 0.9999935461405631
$data[5,9] = 0.9999953532520527
$data[5,10] = 0.9999802549430595
$data[5,11] = 0.999980181993475
$data[5,12] = 1.000034294730554
$data[5,13] = 0.999998628738425
$data[5,14] = 1.000013920435558
$data[5,15] = 0.9999924798234417
$data[5,16] = 0.9999969345391376
$data[5,17] = 0.9999924798234417
$data[5,18] = 0.9999927464027221
$data[5,19] = 0.9999902335208727
$data[5,20] = 0.9999964479696453
$data[6,0] = 1.000047493195834
$data[6,1] = 1.000152453616179
$data[6,2] = 1.000152453616179
$data[6,3] = 0.9999884561728483
$data[6,4] = 0.9998511278963021
$data[6,5] = 0.9998511278963021
$data[6,6] = 0.9998511278963021
$data[6,7] = 0.999946200109595
$data[6,8] = 1.00003435517359
$data[6,9] = 0.9999471284895873
$data[6,10] = 1.000103728428504
$data[6,11] = 0.9998511278963021
$data[6,12] = 0.9999884561728483
$data[6,13] = 1.000070454894514
$data[6,14] = 1.000011405673219
$data[6,15] = 0.9999973458951098
$data[6,16] = 1.000058421654206
$data[6,17] = 0.9999973458951098
$data[6,18] = 1.00000659821473
$data[6,19] = 0.9999755041510443
$data[6,20] = 1.000008867885305
$data[7,0] = 1.000071521503119
$data[7,1] = 1.000141760205109
$data[7,2] = 1.000141760205109
$data[7,3] = 1.000032012198073
$data[7,4] = 0.9997855729726347
$data[7,5] = 0.9997855729726347
$data[7,6] = 0.9997855729726347
$data[7,7] = 0.9999544598172817
$data[7,8] = 1.000034515336138
$data[7,9] = 0.9999270877467316
$data[7,10] = 1.000103855032913
$data[7,11] = 0.9997855729726347
$data[7,12] = 1.000032012198073
$data[7,13] = 1.000086886201591
$data[7,14] = 1.000033263767106
$data[7,15] = 0.9999864484586057
$data[7,16] = 1.00006942924644
$data[7,17] = 0.9999864484586057
$data[7,18] = 0.9999984651779888
$data[7,19] = 0.9999558867369179
$data[7,20] = 1.0000063481015
$data[8,0] = 1.000198545047889
$data[8,1] = 1.000702043745861
$data[8,2] = 1.000702043745861
$data[8,3] = 0.9999153296641975
$data[8,4] = 0.9993704709480303
$data[8,5] = 0.9993704709480303
$data[8,6] = 0.9993704709480303
$data[8,7] = 0.9997489440449054
$data[8,8] = 1.000156311448642
$data[8,9] = 0.9997740197104076
$data[8,10] = 1.000472206316395
$data[8,11] = 0.9993704709480303
$data[8,12] = 0.9999153296641975
$data[8,13] = 1.000308686705029
$data[8,14] = 1.00003582055642
$data[8,15] = 0.9999959481193629
$data[8,16] = 1.0002578949529
$data[8,17] = 0.9999959481193629
$data[8,18] = 1.000036038951683
$data[8,19] = 0.9999029253509523
$data[8,20] = 1.000042233865791
$data[9,0] = 0.9992346878604407
$data[9,1] = 0.9973719629308202
$data[9,2] = 0.9973719629308202
$data[9,3] = 1.000282481962374
$data[9,4] = 1.002417857854792
$data[9,5] = 1.002417857854792
$data[9,6] = 1.002417857854792
$data[9,7] = 1.000936190435607
$data[9,8] = 0.9994127906670638
$data[9,9] = 1.000865077468418
$data[9,10] = 0.9982263831254611
$data[9,11] = 1.002417857854792
$data[9,12] = 1.000282481962374
$data[9,13] = 0.9988272224465973
$data[9,14] = 0.9998476363147191
$data[9,15] = 1.000024100915995
$data[9,16] = 0.9990224118534194
$data[9,17] = 1.000024100915995
$data[9,18] = 0.9998712733537625
$data[9,19] = 1.000380590253968
$data[9,20] = 0.9998434290381222
$data[10,0] = 0.9999319837421722
$data[10,1] = 0.9990817833146899
$data[10,2] = 0.9990817833146899
$data[10,3] = 1.000410226801606
$data[10,4] = 1.000290957719287
$data[10,5] = 1.000290957719287
$data[10,6] = 1.000290957719287
$data[10,7] = 1.000359801015239
$data[10,8] = 0.9998135526204949
$data[10,9] = 1.000129187227138
$data[10,10] = 0.9994342523920381
$data[10,11] = 1.000290957719287
$data[10,12] = 1.000410226801606
$data[10,13] = 0.999746005058148
$data[10,14] = 1.000111889711051
$data[10,15] = 0.9999276559451943
$data[10,16] = 0.9997685209122636
$data[10,17] = 0.9999276559451943
$data[10,18] = 0.9998991301140194
$data[10,19] = 0.9999774956350729
$data[10,20] = 0.9999314681040831
$data[11,0] = 1.053822836603556
$data[11,1] = 1.277954329490422
$data[11,2] = 1.277954329490422
$data[11,3] = 0.9277488733131457
$data[11,4] = 0.8196082955855434
$data[11,5] = 0.8196082955855434
$data[11,6] = 0.8196082955855434
$data[11,7] = 0.8965358238763459
$data[11,8] = 1.059559687412855
$data[11,9] = 0.9320461621327375
$data[11,10] = 1.180249211616174
$data[11,11] = 0.8196082955855434
$data[11,12] = 0.9277488733131457
$data[11,13] = 1.102851601401784
$data[11,14] = 0.9936542803630006
$data[11,15] = 1.008437166129704
$data[11,16] = 1.088420963405474
$data[11,17] = 1.008437166129704
$data[11,18] = 1.021217796450492
$data[11,19] = 0.980895896277502
$data[11,20] = 1.018440652503847
$data[12,0] = 1.015300701163529
$data[12,1] = 1.072680463717056
$data[12,2] = 1.072680463717056
$data[12,3] = 0.9830245918534591
$data[12,4] = 0.9494224432153225
$data[12,5] = 0.9494224432153225
$data[12,6] = 0.9494224432153225
$data[12,7] = 0.9731470645764246
$data[12,8] = 1.015689095667364
$data[12,9] = 0.9811661573293712
$data[12,10] = 1.04746410601897
$data[12,11] = 0.9494224432153225
$data[12,12] = 0.9830245918534591
$data[12,13] = 1.027852527785258
$data[12,14] = 0.9993568437604115
$data[12,15] = 1.001709166261946
$data[12,16] = 1.023798050412626
$data[12,17] = 1.001709166261946
$data[12,18] = 1.0052041486133
$data[12,19] = 0.9940478075337047
$data[12,20] = 1.004736827942687
$data[13,0] = 0.9922280446287014
$data[13,1] = 0.9689097110013526
$data[13,2] = 0.9689097110013526
$data[13,3] = 1.005344596364361
$data[13,4] = 1.025043251916226
$data[13,5] = 1.025043251916226
$data[13,6] = 1.025043251916226
$data[13,7] = 1.011285572788407
$data[13,8] = 0.9931735139939232
$data[13,9] = 1.009121441050548
$data[13,10] = 0.9793645282457094
$data[13,11] = 1.025043251916226
$data[13,12] = 1.005344596364361
$data[13,13] = 0.9871271536828568
$data[13,14] = 0.999259055179142
$data[13,15] = 0.9997658530939798
$data[13,16] = 0.9891426071198789
$data[13,17] = 0.9997658530939798
$data[13,18] = 0.9981177683189656
$data[13,19] = 1.003502865038418
$data[13,20] = 0.9980588324986535
$data[14,0] = 0.9940362334840839
$data[14,1] = 0.9741518358116981
$data[14,2] = 0.9741518358116981
$data[14,3] = 1.005221200660163
$data[14,4] = 1.019438040978662
$data[14,5] = 1.019438040978662
$data[14,6] = 1.019438040978662
$data[14,7] = 1.009464357108274
$data[14,8] = 0.9943712719046848
$data[14,9] = 1.00715139285431
$data[14,10] = 0.9829785330659314
$data[14,11] = 1.019438040978662
$data[14,12] = 1.005221200660163
$data[14,13] = 0.9896865182359307
$data[14,14] = 0.999796236282424
$data[14,15] = 0.9996036924835078
$data[14,16] = 0.991248102792182
$data[14,17] = 0.9996036924835078
$data[14,18] = 0.998295587338802
$data[14,19] = 1.002524078066774
$data[14,20] = 0.998351608233476
$data[15,0] = 0.9977043880564197
$data[15,1] = 0.9849810648624727
$data[15,2] = 0.9849810648624727
$data[15,3] = 1.004861258239289
$data[15,4] = 1.008045461036991
$data[15,5] = 1.008045461036991
$data[15,6] = 1.008045461036991
$data[15,7] = 1.005691046831508
$data[15,8] = 0.9968392999292747
$data[15,9] = 1.00313999994478
$data[15,10] = 0.9904261919714085
$data[15,11] = 1.008045461036991
$data[15,12] = 1.004861258239289
$data[15,13] = 0.9949211615508808
$data[15,14] = 1.000850279084282
$data[15,15] = 0.999295928046251
$data[15,16] = 0.9955605410103455
$data[15,17] = 0.9992959280462509
$data[15,18] = 0.9986817710170068
$data[15,19] = 1.000554509021004
$data[15,20] = 0.9989610888590179
$data[16,0] = 0.9975032763085385
$data[16,1] = 0.9825261832867465
$data[16,2] = 0.9825261832867465
$data[16,3] = 1.005927892849819
$data[16,4] = 1.008876871757501
$data[16,5] = 1.008876871757501
$data[16,6] = 1.008876871757501
$data[16,7] = 1.006649818581047
$data[16,8] = 0.9963390256845401
$data[16,9] = 1.003502097975585
$data[16,10] = 0.9889084342968603
$data[16,11] = 1.008876871757501
$data[16,12] = 1.005927892849819
$data[16,13] = 0.9942270380682826
$data[16,14] = 1.001133459267179
$data[16,15] = 0.9991103159646887
$data[16,16] = 0.9949310339403684
$data[16,17] = 0.9991103159646887
$data[16,18] = 0.9984174933946516
$data[16,19] = 1.000509369067222
$data[16,20] = 0.9987792000925797
$data[17,0] = 1.000017276893865
$data[17,1] = 0.9999655756815365
$data[17,2] = 0.9999655756815365
$data[17,3] = 1.000046361193685
$data[17,4] = 0.9999558269506368
$data[17,5] = 0.9999558269506368
$data[17,6] = 0.9999558269506368
$data[17,7] = 1.000016741574441
$data[17,8] = 0.9999948721144788
$data[17,9] = 0.9999876331927574
$data[17,10] = 0.9999841562078305
$data[17,11] = 0.9999558269506368
$data[17,12] = 1.000046361193685
$data[17,13] = 1.000005968437611
$data[17,14] = 1.000020616654082
$data[17,15] = 0.9999892546086194
$data[17,16] = 1.000002269663234
$data[17,17] = 0.9999892546086194
$data[17,18] = 0.9999906589850843
$data[17,19] = 0.9999836925781949
$data[17,20] = 0.9999960554761538
$data[18,0] = 0.9999532137378868
$data[18,1] = 0.9998779240781555
$data[18,2] = 0.9998779240781555
$data[18,3] = 0.9999955658476287
$data[18,4] = 1.000143519341251
$data[18,5] = 1.000143519341251
$data[18,6] = 1.000143519341251
$data[18,7] = 1.000041645095209
$data[18,8] = 0.9999716678036198
$data[18,9] = 1.000049939544589
$data[18,10] = 0.9999145757953288
$data[18,11] = 1.000143519341251
$data[18,12] = 0.9999955658476287
$data[18,13] = 0.9999367449628921
$data[18,14] = 0.9999836168256242
$data[18,15] = 1.000005669755679
$data[18,16] = 0.9999483859098013
$data[18,17] = 1.000005669755678
$data[18,18] = 0.9999971692676637
$data[18,19] = 1.000026439282381
$data[18,20] = 0.9999935064054586
$data[19,0] = 0.9999185844993872
$data[19,1] = 0.9994477109054802
$data[19,2] = 0.9994477109054802
$data[19,3] = 1.000183448643837
$data[19,4] = 1.000287516832997
$data[19,5] = 1.000287516832997
$data[19,6] = 1.000287516832997
$data[19,7] = 1.000209768435695
$data[19,8] = 0.9998840541412201
$data[19,9] = 1.000112863738293
$data[19,10] = 0.999648759831062
$data[19,11] = 1.000287516832997
$data[19,12] = 1.000183448643837
$data[19,13] = 0.9998155797746584
$data[19,14] = 1.000033751392528
$data[19,15] = 0.9999728921274378
$data[19,16] = 0.9998384045635124
$data[19,17] = 0.9999728921274378
$data[19,18] = 0.9999506826308834
$data[19,19] = 1.000018049471306
$data[19,20] = 0.9999615883784964
$data[20,0] = 0.9998502573002208
$data[20,1] = 0.9986458786381561
$data[20,2] = 0.9986458786381561
$data[20,3] = 1.000527720428362
$data[20,4] = 1.000566402307763
$data[20,5] = 1.000566402307763
$data[20,6] = 1.000566402307763
$data[20,7] = 1.000522497441397
$data[20,8] = 0.9997204040067069
$data[20,9] = 1.000233419849397
$data[20,10] = 0.9991523043883592
$data[20,11] = 1.000566402307763
$data[20,12] = 1.000527720428362
$data[20,13] = 0.9995867995332592
$data[20,14] = 1.000124062217535
$data[20,15] = 0.9999133337914273
$data[20,16] = 0.9996313343577418
$data[20,17] = 0.9999133337914273
$data[20,18] = 0.9998651013452473
$data[20,19] = 1.00000536153775
$data[20,20] = 0.9999023605450453
$ws.Range("C3:W23").Value = $data

Write-Host "edit applied"
